# Postprocess update: refresh mapping-time/total-time figures for the
# BFAST/BWA/BOWTIE rows and append the (now superseded) original BFAST
# totals as a trailing block at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated figures for the existing rows (mouse ch19 / 100MB block) ---
$ws.Range("J6").Value = 44
$ws.Range("L6").Value = 152.78

$ws.Range("J7").Value = 84
$ws.Range("L7").Value = 217.66

# --- New figures for the 200MB block (rows 10-11), previously blank ---
$ws.Range("J10").Value = 33
$ws.Range("L10").Value = 142.34

$ws.Range("J11").Value = 52
$ws.Range("L11").Value = 202.22

# --- New figures for the 400MB block (rows 14-15), previously blank ---
$ws.Range("J14").Value = 30
$ws.Range("L14").Value = 137.33

$ws.Range("J15").Value = 60
$ws.Range("L15").Value = 177.45

# --- Trailing rows 41-42: previous (pre-postprocess) BFAST totals ---
$ws.Range("J41").Value = 68
$ws.Range("L41").Value = 181.7

$ws.Range("J42").Value = 35
$ws.Range("L42").Value = 136.41

# --- Update the active selection left by the editor ---
$ws.Range("F29").Select()
